$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (Price, Volume) updates for the periodic cryptos list refresh.
# Price values that look numeric are prefixed with a leading apostrophe so
# Excel stores them as text (matching the original inline-string cells)
# instead of silently converting them to numbers.
$updates = @(
    @{Row=2; D="25.955.61"; E="  +0.74%  "}
    @{Row=3; D="1.747.14"; E="  -0.12%  "}
    @{Row=4; D="'0.9999"; E="  -0.16%  "}
    @{Row=5; D="'233.77"; E="  -1.17%  "}
    @{Row=6; D="'0.9993"; E="  -0.19%  "}
    @{Row=7; D="'0.5169"; E="  +2.22%  "}
    @{Row=8; D="'0.2826"; E="  +8.47%  "}
    @{Row=9; D="'39.80"; E="  -1.68%  "}
    @{Row=10; D="'0.06121"; E="  -0.86%  "}
    @{Row=11; D="1.752.29"; E="  +0.15%  "}
    @{Row=12; D="'0.07021"; E="  +1.28%  "}
    @{Row=13; D="'15.43"; E="  +0.46%  "}
    @{Row=14; D="'0.6422"; E="  +6.06%  "}
    @{Row=15; D="'4.515"; E="  +1.41%  "}
    @{Row=16; D="'77.12"; E="  -1.62%  "}
    @{Row=17; D="'0.9981"; E="  -0.29%  "}
    @{Row=18; D="'0.9981"; E="  -0.25%  "}
    @{Row=19; D="25.969.43"; E="  +0.72%  "}
    @{Row=20; D="'11.50"; E="  -1.26%  "}
    @{Row=21; D="'0.000006608"; E="  -1.12%  "}
    @{Row=22; D="1.967.44"; E="  -0.44%  "}
    @{Row=23; D="'4.146"; E="  +2.42%  "}
    @{Row=24; D="'8.566"; E="  +4.53%  "}
    @{Row=25; D="'5.149"; E="  -0.13%  "}
    @{Row=26; D="'140.19"; E="  +2.03%  "}
    @{Row=27; D="'1.492"; E="  +2.35%  "}
    @{Row=28; D="'1.840"; E="  +2.22%  "}
    @{Row=29; D="'15.07"; E="  -0.08%  "}
    @{Row=30; D="'103.13"; E="  +1.02%  "}
    @{Row=31; D="'0.08296"; E="  +0.36%  "}
    @{Row=32; D="'3.645"; E="  -1.42%  "}
    @{Row=33; D="'3.428"; E="  +1.00%  "}
    @{Row=34; D="'0.04412"; E="  +1.26%  "}
    @{Row=35; D="'2.605"; E="  -1.80%  "}
    @{Row=36; D="'0.9839"; E="  -1.29%  "}
    @{Row=37; D="'0.6096"; E="  +1.57%  "}
    @{Row=38; D="'2.691"; E="  -0.30%  "}
    @{Row=39; D="'0.01576"; E="  +1.72%  "}
    @{Row=40; D="'1.934"; E="  -1.00%  "}
    @{Row=41; D="'0.9977"; E="  -0.34%  "}
    @{Row=42; D="'100.68"; E="  -1.93%  "}
    @{Row=43; D="'0.3864"; E="  +1.63%  "}
    @{Row=44; D="'0.7348"; E="  -2.93%  "}
    @{Row=45; D="'4.983"; E="  +2.80%  "}
    @{Row=46; D="'0.05469"; E="  -0.46%  "}
    @{Row=47; D="'6.360"; E="  +7.55%  "}
    @{Row=48; D="'0.1118"; E="  +3.70%  "}
    @{Row=49; D="'52.68"; E="  +1.21%  "}
    @{Row=50; D="'29.93"; E="  -0.68%  "}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

# Row 51: coin listing swapped from Decentraland to EnergySwap
$ws.Cells.Item(51, 2).Value = "EnergySwap"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(51, 4).Value = "'7.530"
$ws.Cells.Item(51, 5).Value = "  +1.03%  "
